# Daily attendance processing - 2025-10-15 06:57:26
# Applies the latest attendance-recording results to the
# "Session Analysis Results" sheet: updates "Recorded By" lists,
# attendance counts, two sessions that moved from Pending/Not Recorded
# to Recorded (rows 5 and 33), and the derived summary statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a percentage-looking (or otherwise numeric-looking) piece of
# text into a cell while keeping it as literal text (matching the workbook's
# existing convention of storing percentages as plain strings, not numeric
# percentage values) and keeping the cell's original visual style.
function Set-TextValue($cellAddr, $val, $donorAddr) {
    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $val
    $ws.Range($donorAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Row 2 - ANATOMY session 1 (Year 3 / C1): recorder list updated
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 3 - ANATOMY session 2 (Year 3 / C1): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H3").Value = "108/221"

# ---------------------------------------------------------------------
# Row 4 - ANATOMY session 3 (Year 3 / C1): recorder list updated
# ---------------------------------------------------------------------
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 5 - ANATOMY session 4 (Year 3 / C1): moved from Pending -> Recorded
# Copy the "Recorded" look (style) from row 2, then fill in the new data.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H5").Value = "3/221"
$ws.Range("I5").Value = "Recorded"

# ---------------------------------------------------------------------
# Class Statistics block (rows 6-10) for Year 3 / C1
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 40
Set-TextValue "L9" "26.8%" "L8"
Set-TextValue "L10" "41.5%" "L8"

# ---------------------------------------------------------------------
# Row 13 - HISTOLOGY session 1 (Year 3 / C1): recorder list updated
# ---------------------------------------------------------------------
$ws.Range("G13").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Group Statistics block: row 15 (Year 3 / C1), row 16 (Year 3 / C2)
# ---------------------------------------------------------------------
$ws.Range("O15").Value = 8
$ws.Range("Q15").Value = 20
Set-TextValue "R15" "28.6%" "Q15"
Set-TextValue "S15" "46.4%" "Q15"

$ws.Range("O16").Value = 7
$ws.Range("P16").Value = 1
Set-TextValue "R16" "25.0%" "Q16"
Set-TextValue "S16" "35.9%" "Q16"

# ---------------------------------------------------------------------
# Row 24 - PHYSIOLOGY session 1 (Year 3 / C1): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G24").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("H24").Value = "150/221"

# ---------------------------------------------------------------------
# Row 25 - PHYSIOLOGY session 2 (Year 3 / C1): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G25").Value = "abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("H25").Value = "89/221"

# ---------------------------------------------------------------------
# Row 31 - ANATOMY session 2 (Year 3 / C2): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G31").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H31").Value = "63/246"

# ---------------------------------------------------------------------
# Row 32 - ANATOMY session 3 (Year 3 / C2): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G32").Value = "hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("H32").Value = "122/246"

# ---------------------------------------------------------------------
# Row 33 - ANATOMY session 4 (Year 3 / C2): moved from Not Recorded -> Recorded
# Copy the "Recorded" look (style) from row 31 (PasteSpecial formats-only
# leaves this row's own A/D/E/F values - Year 3 / C2 / ANATOMY / 4 /
# 15/10/2025 / 08:00:00 - untouched), then fill in the new data.
# ---------------------------------------------------------------------
$ws.Range("A31:I31").Copy()
$ws.Range("A33:I33").PasteSpecial(-4122) | Out-Null
$ws.Range("G33").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H33").Value = "46/246"
$ws.Range("I33").Value = "Recorded"

# ---------------------------------------------------------------------
# Row 41 - HISTOLOGY session 1 (Year 3 / C2): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G41").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("H41").Value = "71/246"

# ---------------------------------------------------------------------
# Row 52 - PHYSIOLOGY session 1 (Year 3 / C2): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G52").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("H52").Value = "101/246"

# ---------------------------------------------------------------------
# Row 53 - PHYSIOLOGY session 2 (Year 3 / C2): recorder list + students
# ---------------------------------------------------------------------
$ws.Range("G53").Value = "abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("H53").Value = "70/246"
